# Finance Change Management Plan -> Banking Implementation Project
# Reverts the "AI/ML" / "FINANCE" branded template text back to the
# generic "Banking" template text, and restores a few blank spacer rows
# that exist in the canonical template.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Change Management Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A2").Value = "Banking Implementation Project"
$ws1.Range("B6").Value = "Enterprise Banking Implementation"
$ws1.Range("A15").Value = "1. Achieve 95% user adoption of new Banking systems within 6 months of go-live"
$ws1.Range("A17").Value = "3. Build organizational capability and confidence in Banking technologies"
$ws1.Range("A20").Value = "6. Create positive stakeholder sentiment and enthusiasm for Banking transformation"

# Restore blank spacer row 4, which the original template keeps as an
# explicit (empty) row between the "Generated:" line and the
# "PROJECT INFORMATION" section header.
$ws1.Rows.Item(4).OutlineLevel = 1
$ws1.Rows.Item(4).OutlineLevel = 0

# Restore blank spacer row 13, between the project info block and the
# "CHANGE MANAGEMENT OBJECTIVES" section header.
$ws1.Rows.Item(13).OutlineLevel = 1
$ws1.Rows.Item(13).OutlineLevel = 0

# Restore blank spacer row 21, between the objectives list and the
# "CHANGE MANAGEMENT STRATEGY" section header.
$ws1.Rows.Item(21).OutlineLevel = 1
$ws1.Rows.Item(21).OutlineLevel = 0

# ---------------------------------------------------------------------
# Sheet 2: "Change Impact Assessment"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Restore blank spacer row 2, between the title and the table header row.
$ws2.Rows.Item(2).OutlineLevel = 1
$ws2.Rows.Item(2).OutlineLevel = 0

$ws2.Range("G4").Value = "Banking automation"

# ---------------------------------------------------------------------
# Sheet 3: "Change Activities"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

# Restore blank spacer row 2, between the title and the table header row.
$ws3.Rows.Item(2).OutlineLevel = 1
$ws3.Rows.Item(2).OutlineLevel = 0
